# "fixed export and fixing maps"
#
# The sheet still carried the generic default name "1" and the table
# included an obsolete census-results subtitle plus two stale area
# columns (1989 and 2002) that shouldn't be exported any more - only
# the 2014 figure is current. Clean that up:

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the sheet its real (Georgian) name instead of the default "1".
$ws.Name = "თელავი"

# Drop the "(census results)" subtitle row entirely - everything below
# shifts up one row.
$ws.Rows(2).Delete()

# Drop the 1989 and 2002 columns - only the 2014 figures are still
# published, and they shift left into column B.
$ws.Columns("B:C").Delete()

# Leave the selection where it lands after removing the subtitle row.
[void]$ws.Range("A2").Select()
